$wb = $excel.ActiveWorkbook

$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

# --- Sheet "About" ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: " + $newVersion

$wsAbout.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for No. 4 Coal Mine (AL), United States, M1397, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet "Boundaries and methane sources" ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 27; $row++) {
    $wsData.Range("S" + $row).Value = $newVersion
}
